$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the missing "received_final" (column K) flags for rows 44, 48-56 ---
$ws.Range("K44").Value = 1
$ws.Range("K48").Value = 1
$ws.Range("K49").Value = 1
$ws.Range("K50").Value = 1
$ws.Range("K51").Value = 1
$ws.Range("K52").Value = 1
$ws.Range("K53").Value = 1
$ws.Range("K54").Value = 1
$ws.Range("K55").Value = 1
$ws.Range("K56").Value = 1

# --- Add new row 57 for survey_round 44 (week 45, panel F, wave 13) ---
$ws.Range("A57").Value = 3
$ws.Range("B57").Value = 0
$ws.Range("C57").Value = "uk"
$ws.Range("D57").Value = 45
$ws.Range("E57").Value = "F"
$ws.Range("F57").Formula = "=F55+1"

# Copy the date style from the row above so the new date cell reuses the
# existing date number-format style instead of creating a new one.
$ws.Range("G56").Copy()
$ws.Range("G57").PasteSpecial(-4122)
$ws.Range("G57").Value = 44232

$ws.Range("H57").Value = "20-100590_PFW13_Final_ICUO"
$ws.Range("I57").Formula = '=C57&"_"&"wk"&TEXT(D57,"00")&"_"&YEAR(G57)&TEXT(G57,"MM")&TEXT(G57,"DD")&"_p"&E57&"_wv"&TEXT(F57,"00")&""'
$ws.Range("J57").Value = 1
$ws.Range("K57").Value = 1

# --- Match the author's final on-screen selection state ---
$ws.Range("K48:K57").Select() | Out-Null
